$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 12286
$ws.Range("F3").Value = 6906
$ws.Range("F7").Value = 262
$ws.Range("F11").Value = 119
$ws.Range("F13").Value = 967
$ws.Range("F14").Value = 3697
$ws.Range("F18").Value = 214
$ws.Range("F22").Value = 283
$ws.Range("F23").Value = 71
$ws.Range("F24").Value = 332
$ws.Range("F25").Value = 5112
$ws.Range("F27").Value = 1326
$ws.Range("F28").Value = 265
$ws.Range("F29").Value = 789
$ws.Range("F30").Value = 1276

# Sheet 2: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 3720
$ws.Range("F6").Value = 12
$ws.Range("F7").Value = 26

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9174
$ws.Range("F4").Value = 1909

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9174
$ws.Range("F4").Value = 1909
$ws.Range("F5").Value = 12286
$ws.Range("F6").Value = 6906
$ws.Range("F8").Value = 3720
$ws.Range("F12").Value = 262
$ws.Range("F16").Value = 119
$ws.Range("F18").Value = 967
$ws.Range("F19").Value = 3697
$ws.Range("F23").Value = 214
$ws.Range("F27").Value = 283
$ws.Range("F29").Value = 12
$ws.Range("F32").Value = 332
$ws.Range("F33").Value = 5112
$ws.Range("F35").Value = 1326
$ws.Range("F38").Value = 265
$ws.Range("F40").Value = 789
$ws.Range("F41").Value = 1276

$wb.Save()
